$d = $word.ActiveDocument

# 1. Update the talk date in the title-page "Date" paragraph.
$null = $d.Content.Find.Execute("10/9/2014", $false, $false, $false, $false, $false, `
                                 $true, 1, $false, "15/9/2014", 2)

# 2. Append " http://github.io/tverbeiren/ReproducibleDataAnalysis/" as a live
#    hyperlink right after the last paragraph's existing text ("You can find
#    everything I showed here at:").
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$tail = $lastPara.Range
$url = "http://github.io/tverbeiren/ReproducibleDataAnalysis/"

# Type the separating space plus the raw URL text first …
$tail.InsertAfter(" " + $url)

# … then carve out just the URL portion (Content.End sits one position past
# the last real character, at the phantom end-of-story mark) and turn it
# into a hyperlink, matching the style used by the other links in this doc.
$storyEnd = $d.Content.End - 1
$urlStart = $storyEnd - $url.Length
$urlRange = $d.Range($urlStart, $storyEnd)

$null = $d.Hyperlinks.Add($urlRange, $url, $null, $null, $url)
$urlRange.Style = "Link"
